$d = $word.ActiveDocument

$d.Content.Find.Execute("60×67=", $true, $false, $false, $false, $false, $true, 1, $false, "58×59=", 2) | Out-Null
$d.Content.Find.Execute("86×91=", $true, $false, $false, $false, $false, $true, 1, $false, "44×65=", 2) | Out-Null
$d.Content.Find.Execute("74×40=", $true, $false, $false, $false, $false, $true, 1, $false, "50×53=", 2) | Out-Null
$d.Content.Find.Execute("96×81=", $true, $false, $false, $false, $false, $true, 1, $false, "35×65=", 2) | Out-Null
$d.Content.Find.Execute("92×11=", $true, $false, $false, $false, $false, $true, 1, $false, "41×32=", 2) | Out-Null
$d.Content.Find.Execute("61×91=", $true, $false, $false, $false, $false, $true, 1, $false, "77×34=", 2) | Out-Null
$d.Content.Find.Execute("94×40=", $true, $false, $false, $false, $false, $true, 1, $false, "17×41=", 2) | Out-Null
$d.Content.Find.Execute("50×11=", $true, $false, $false, $false, $false, $true, 1, $false, "99×76=", 2) | Out-Null
$d.Content.Find.Execute("73×61=", $true, $false, $false, $false, $false, $true, 1, $false, "90×76=", 2) | Out-Null
$d.Content.Find.Execute("49×55=", $true, $false, $false, $false, $false, $true, 1, $false, "67×94=", 2) | Out-Null
$d.Content.Find.Execute("44×15=", $true, $false, $false, $false, $false, $true, 1, $false, "56×54=", 2) | Out-Null
$d.Content.Find.Execute("31×46=", $true, $false, $false, $false, $false, $true, 1, $false, "50×86=", 2) | Out-Null
$d.Content.Find.Execute("74×78=", $true, $false, $false, $false, $false, $true, 1, $false, "11×59=", 2) | Out-Null
$d.Content.Find.Execute("31×71=", $true, $false, $false, $false, $false, $true, 1, $false, "59×50=", 2) | Out-Null
$d.Content.Find.Execute("74×35=", $true, $false, $false, $false, $false, $true, 1, $false, "30×60=", 2) | Out-Null
$d.Content.Find.Execute("49×92=", $true, $false, $false, $false, $false, $true, 1, $false, "42×53=", 2) | Out-Null
$d.Content.Find.Execute("44×91=", $true, $false, $false, $false, $false, $true, 1, $false, "70×23=", 2) | Out-Null
$d.Content.Find.Execute("88×29=", $true, $false, $false, $false, $false, $true, 1, $false, "51×50=", 2) | Out-Null
$d.Content.Find.Execute("87×35=", $true, $false, $false, $false, $false, $true, 1, $false, "41×34=", 2) | Out-Null
$d.Content.Find.Execute("64×90=", $true, $false, $false, $false, $false, $true, 1, $false, "87×89=", 2) | Out-Null
$d.Content.Find.Execute("34×71=", $true, $false, $false, $false, $false, $true, 1, $false, "50×79=", 2) | Out-Null
$d.Content.Find.Execute("44×94=", $true, $false, $false, $false, $false, $true, 1, $false, "22×85=", 2) | Out-Null
$d.Content.Find.Execute("32×93=", $true, $false, $false, $false, $false, $true, 1, $false, "28×30=", 2) | Out-Null
$d.Content.Find.Execute("98×30=", $true, $false, $false, $false, $false, $true, 1, $false, "21×15=", 2) | Out-Null
$d.Content.Find.Execute("84×93=", $true, $false, $false, $false, $false, $true, 1, $false, "95×91=", 2) | Out-Null
